$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = 130746519
$ws.Range("M10").Value = "färska gnagspår"
$ws.Range("Q10").Value = 447826
$ws.Range("R10").Value = 6784623

# Row 12
$ws.Range("A12").Value = 130746562
$ws.Range("B12").Value = 79243
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = "Garnlav"
$ws.Range("G12").Value = "Alectoria sarmentosa"
$ws.Range("H12").Value = "(Ach.) Ach."
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("L12").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("Q12").Value = 447730
$ws.Range("R12").Value = 6784717
$ws.Range("AF12").Value = ""

# Row 13
$ws.Range("A13").Value = 130746524
$ws.Range("B13").Value = 8451
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 106545
$ws.Range("F13").Value = "Mindre märgborre"
$ws.Range("G13").Value = "Tomicus minor"
$ws.Range("H13").Value = "(Hartig, 1834)"
$ws.Range("J13").NumberFormat = "General"
$ws.Range("K13").NumberFormat = "General"
$ws.Range("L13").NumberFormat = "General"
$ws.Range("M13").Value = "äldre gnagspår"
$ws.Range("N13").NumberFormat = "General"
$ws.Range("Q13").Value = 447932
$ws.Range("R13").Value = 6784551
$ws.Range("AF13").NumberFormat = "General"

# Row 16
$ws.Range("A16").Value = 130746499
$ws.Range("B16").Value = 57881
$ws.Range("E16").Value = 100049
$ws.Range("F16").Value = "Spillkråka"
$ws.Range("G16").Value = "Dryocopus martius"
$ws.Range("H16").Value = "(Linnaeus, 1758)"
$ws.Range("J16").Value = ""
$ws.Range("L16").NumberFormat = "General"
$ws.Range("M16").Value = "äldre spår"
$ws.Range("Q16").Value = 447930
$ws.Range("R16").Value = 6784568
$ws.Range("AF16").Value = ""

# Row 17
$ws.Range("A17").Value = 130746554
$ws.Range("B17").Value = 79243
$ws.Range("E17").Value = 6425
$ws.Range("F17").Value = "Garnlav"
$ws.Range("G17").Value = "Alectoria sarmentosa"
$ws.Range("H17").Value = "(Ach.) Ach."
$ws.Range("J17").NumberFormat = "General"
$ws.Range("L17").Value = ""
$ws.Range("M17").Value = ""
$ws.Range("Q17").Value = 447915
$ws.Range("R17").Value = 6784490
$ws.Range("AF17").NumberFormat = "General"

# Row 18
$ws.Range("A18").Value = 130746500
$ws.Range("Q18").Value = 447888
$ws.Range("R18").Value = 6784627

# Row 25
$ws.Range("A25").Value = 130746515
$ws.Range("B25").Value = 8451
$ws.Range("D25").Value = "LC"
$ws.Range("E25").Value = 106545
$ws.Range("F25").Value = "Mindre märgborre"
$ws.Range("G25").Value = "Tomicus minor"
$ws.Range("H25").Value = "(Hartig, 1834)"
$ws.Range("J25").NumberFormat = "General"
$ws.Range("K25").NumberFormat = "General"
$ws.Range("L25").NumberFormat = "General"
$ws.Range("M25").Value = "färska gnagspår"
$ws.Range("N25").NumberFormat = "General"
$ws.Range("Q25").Value = 447716
$ws.Range("R25").Value = 6784496
$ws.Range("AF25").NumberFormat = "General"

# Row 26
$ws.Range("A26").Value = 130746518
$ws.Range("B26").Value = 8451
$ws.Range("D26").Value = "LC"
$ws.Range("E26").Value = 106545
$ws.Range("F26").Value = "Mindre märgborre"
$ws.Range("G26").Value = "Tomicus minor"
$ws.Range("H26").Value = "(Hartig, 1834)"
$ws.Range("J26").NumberFormat = "General"
$ws.Range("K26").NumberFormat = "General"
$ws.Range("L26").NumberFormat = "General"
$ws.Range("M26").Value = "äldre gnagspår"
$ws.Range("N26").NumberFormat = "General"
$ws.Range("Q26").Value = 447815
$ws.Range("R26").Value = 6784612
$ws.Range("AF26").NumberFormat = "General"

# Row 27
$ws.Range("A27").Value = 130746569
$ws.Range("B27").Value = 79243
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("J27").Value = ""
$ws.Range("K27").Value = ""
$ws.Range("L27").Value = ""
$ws.Range("M27").Value = ""
$ws.Range("N27").Value = ""
$ws.Range("Q27").Value = 447856
$ws.Range("R27").Value = 6784518
$ws.Range("AF27").Value = ""

# Row 28
$ws.Range("A28").Value = 130746564
$ws.Range("B28").Value = 79243
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = "Garnlav"
$ws.Range("G28").Value = "Alectoria sarmentosa"
$ws.Range("H28").Value = "(Ach.) Ach."
$ws.Range("J28").Value = ""
$ws.Range("K28").Value = ""
$ws.Range("L28").Value = ""
$ws.Range("M28").Value = ""
$ws.Range("N28").Value = ""
$ws.Range("Q28").Value = 447866
$ws.Range("R28").Value = 6784648
$ws.Range("AF28").Value = ""

# Row 29
$ws.Range("A29").Value = 130746496
$ws.Range("B29").Value = 57881
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 100049
$ws.Range("F29").Value = "Spillkråka"
$ws.Range("G29").Value = "Dryocopus martius"
$ws.Range("H29").Value = "(Linnaeus, 1758)"
$ws.Range("J29").Value = ""
$ws.Range("M29").Value = "färska spår"
$ws.Range("Q29").Value = 447817
$ws.Range("R29").Value = 6784636
$ws.Range("AF29").Value = ""

# Row 30
$ws.Range("A30").Value = 130746521
$ws.Range("B30").Value = 8451
$ws.Range("D30").Value = "LC"
$ws.Range("E30").Value = 106545
$ws.Range("F30").Value = "Mindre märgborre"
$ws.Range("G30").Value = "Tomicus minor"
$ws.Range("H30").Value = "(Hartig, 1834)"
$ws.Range("J30").NumberFormat = "General"
$ws.Range("M30").Value = "äldre gnagspår"
$ws.Range("Q30").Value = 447834
$ws.Range("R30").Value = 6784604
$ws.Range("AF30").NumberFormat = "General"

# Row 32
$ws.Range("A32").Value = 130746497
$ws.Range("B32").Value = 57881
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 100049
$ws.Range("F32").Value = "Spillkråka"
$ws.Range("G32").Value = "Dryocopus martius"
$ws.Range("H32").Value = "(Linnaeus, 1758)"
$ws.Range("J32").Value = ""
$ws.Range("M32").Value = "färska spår"
$ws.Range("Q32").Value = 447838
$ws.Range("R32").Value = 6784644
$ws.Range("AF32").Value = ""

# Row 34
$ws.Range("A34").Value = 130746523
$ws.Range("B34").Value = 8451
$ws.Range("D34").Value = "LC"
$ws.Range("E34").Value = 106545
$ws.Range("F34").Value = "Mindre märgborre"
$ws.Range("G34").Value = "Tomicus minor"
$ws.Range("H34").Value = "(Hartig, 1834)"
$ws.Range("J34").NumberFormat = "General"
$ws.Range("K34").NumberFormat = "General"
$ws.Range("L34").NumberFormat = "General"
$ws.Range("M34").Value = "äldre gnagspår"
$ws.Range("N34").NumberFormat = "General"
$ws.Range("Q34").Value = 447940
$ws.Range("R34").Value = 6784589
$ws.Range("AF34").NumberFormat = "General"

# Row 35
$ws.Range("A35").Value = 130746555
$ws.Range("B35").Value = 79243
$ws.Range("E35").Value = 6425
$ws.Range("F35").Value = "Garnlav"
$ws.Range("G35").Value = "Alectoria sarmentosa"
$ws.Range("H35").Value = "(Ach.) Ach."
$ws.Range("K35").Value = ""
$ws.Range("L35").Value = ""
$ws.Range("M35").Value = ""
$ws.Range("N35").Value = ""
$ws.Range("Q35").Value = 447906
$ws.Range("R35").Value = 6784505

# Row 36
$ws.Range("A36").Value = 130746565
$ws.Range("B36").Value = 79243
$ws.Range("D36").Value = "NT"
$ws.Range("E36").Value = 6425
$ws.Range("F36").Value = "Garnlav"
$ws.Range("G36").Value = "Alectoria sarmentosa"
$ws.Range("H36").Value = "(Ach.) Ach."
$ws.Range("J36").Value = ""
$ws.Range("K36").Value = ""
$ws.Range("L36").Value = ""
$ws.Range("M36").Value = ""
$ws.Range("N36").Value = ""
$ws.Range("Q36").Value = 447912
$ws.Range("R36").Value = 6784599
$ws.Range("AF36").Value = ""

# Row 38
$ws.Range("A38").Value = 130746561
$ws.Range("Q38").Value = 447711
$ws.Range("R38").Value = 6784677

# Row 39
$ws.Range("A39").Value = 130746506
$ws.Range("B39").Value = 8451
$ws.Range("D39").Value = "LC"
$ws.Range("E39").Value = 106545
$ws.Range("F39").Value = "Mindre märgborre"
$ws.Range("G39").Value = "Tomicus minor"
$ws.Range("H39").Value = "(Hartig, 1834)"
$ws.Range("J39").NumberFormat = "General"
$ws.Range("K39").NumberFormat = "General"
$ws.Range("L39").NumberFormat = "General"
$ws.Range("M39").Value = "äldre gnagspår"
$ws.Range("N39").NumberFormat = "General"
$ws.Range("Q39").Value = 447826
$ws.Range("R39").Value = 6784573
$ws.Range("AF39").NumberFormat = "General"

# Row 43
$ws.Range("A43").Value = 130746511
$ws.Range("B43").Value = 8451
$ws.Range("D43").Value = "LC"
$ws.Range("E43").Value = 106545
$ws.Range("F43").Value = "Mindre märgborre"
$ws.Range("G43").Value = "Tomicus minor"
$ws.Range("H43").Value = "(Hartig, 1834)"
$ws.Range("J43").NumberFormat = "General"
$ws.Range("K43").NumberFormat = "General"
$ws.Range("L43").NumberFormat = "General"
$ws.Range("M43").Value = "äldre gnagspår"
$ws.Range("N43").NumberFormat = "General"
$ws.Range("Q43").Value = 447748
$ws.Range("R43").Value = 6784472
$ws.Range("AF43").NumberFormat = "General"

# Row 44
$ws.Range("A44").Value = 130746566
$ws.Range("B44").Value = 79243
$ws.Range("D44").Value = "NT"
$ws.Range("E44").Value = 6425
$ws.Range("F44").Value = "Garnlav"
$ws.Range("G44").Value = "Alectoria sarmentosa"
$ws.Range("H44").Value = "(Ach.) Ach."
$ws.Range("J44").Value = ""
$ws.Range("K44").Value = ""
$ws.Range("L44").Value = ""
$ws.Range("M44").Value = ""
$ws.Range("N44").Value = ""
$ws.Range("Q44").Value = 447949
$ws.Range("R44").Value = 6784550
$ws.Range("AF44").Value = ""
